# Optuna Attempt (go back with original)
# Update the "Seasonality Index" (column L) values on the "Forecast Comparison" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

$ws.Range("L2").Value  = 0.97
$ws.Range("L3").Value  = 1.08
$ws.Range("L4").Value  = 1
$ws.Range("L5").Value  = 0.96
$ws.Range("L6").Value  = 1.05
$ws.Range("L7").Value  = 1.01
$ws.Range("L8").Value  = 1.07
$ws.Range("L9").Value  = 0.8
$ws.Range("L10").Value = 1.02
$ws.Range("L12").Value = 1.12
$ws.Range("L13").Value = 0.9
$ws.Range("L14").Value = 0.98
$ws.Range("L15").Value = 0.9399999999999999
$ws.Range("L16").Value = 0.93
$ws.Range("L17").Value = 1.08
